# Update "想去人数" (want-to-go count) figures in column F for rows 3-6
# on both the "展览" sheet and the "全部类型" sheet (which mirrors the
# same rows), matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    3 = 128
    4 = 174
    5 = 3289
    6 = 332
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
